# Inserts a new data row at row 161 (shifting the existing rows 161-261
# down to 162-262, so the former last row, 261, becomes row 262) and
# populates the newly inserted row with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 161..261 down by one row, creating a blank row 161.
$ws.Rows.Item(161).Insert()

# Populate the new row 161 with the new record (columns A-T).
$ws.Cells.Item(161, 1).Value = 11
$ws.Cells.Item(161, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(161, 3).Value = "Bíobío"
$ws.Cells.Item(161, 4).Value = 45126
$ws.Cells.Item(161, 5).Value = 8
$ws.Cells.Item(161, 6).Value = "Fruta"
$ws.Cells.Item(161, 7).Value = 100109
$ws.Cells.Item(161, 8).Value = "Uva"
$ws.Cells.Item(161, 9).Value = 100109001
$ws.Cells.Item(161, 10).Value = "Uva"
$ws.Cells.Item(161, 11).Value = "Red Globe"
$ws.Cells.Item(161, 12).Value = "Segunda"
$ws.Cells.Item(161, 13).Value = 100
$ws.Cells.Item(161, 14).Value = 10000
$ws.Cells.Item(161, 15).Value = 11000
$ws.Cells.Item(161, 16).Value = 10500
$ws.Cells.Item(161, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(161, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(161, 19).Value = 1312
$ws.Cells.Item(161, 20).Value = 8
